$d = $word.ActiveDocument

# 1. Add a "_GoBack" bookmark at the very start of the document (before "HENRY HE")
$bookmarkRange = $d.Range(0, 0)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# 2. Merge "JavaScript" + ", " runs (which previously had the bookmark between them)
#    into a single run "JavaScript, " - achieved via Find/Replace.
$d.Content.Find.Execute("JavaScript, ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "JavaScript, ", 2)

# 3. Adjust top margin from 41.85pt (837 twips) to 58.05pt (1161 twips)
$d.PageSetup.TopMargin = 58.05
